# Update template_mahasiswa.xlsx:
#  - remove the "keterangan" column (was column J), shifting prodi_id
#    (formerly column K) left into column J
#  - refresh the sample data row with a new student's details

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "keterangan" column (column J); prodi_id shifts from K -> J.
$ws.Columns.Item(10).Delete()

# Refresh the sample student data in row 2.
$ws.Range("A2").Formula = "2341760026"
$ws.Range("B2").Formula = "' 3507166210050002"
$ws.Range("C2").Formula = "Indi Warda"
$ws.Range("F2").Formula = "Malang"
$ws.Range("G2").Formula = "Malang"

# Update the view so it matches the refreshed sheet.
$ws.Range("K5").Select()
